$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column F, styled like the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Timestamps taken from the source metadata (microsecond-precision, as recorded
# when the panel data was originally processed).
$timestamps = @(
    "2021-10-05 13:39:28.486192",
    "2021-10-05 13:39:28.486206",
    "2021-10-05 13:39:28.486210",
    "2021-10-05 13:39:28.486213",
    "2021-10-05 13:39:28.486216",
    "2021-10-05 13:39:28.486220",
    "2021-10-05 13:39:28.486223",
    "2021-10-05 13:39:28.486226",
    "2021-10-05 13:39:28.486229",
    "2021-10-05 13:39:28.486232",
    "2021-10-05 13:39:28.486235",
    "2021-10-05 13:39:28.486238",
    "2021-10-05 13:39:28.486242",
    "2021-10-05 13:39:28.486245",
    "2021-10-05 13:39:28.486248",
    "2021-10-05 13:39:28.486251",
    "2021-10-05 13:39:28.486254",
    "2021-10-05 13:39:28.486258",
    "2021-10-05 13:39:28.486261",
    "2021-10-05 13:39:28.486264",
    "2021-10-05 13:39:28.486267",
    "2021-10-05 13:39:28.486271",
    "2021-10-05 13:39:28.486274",
    "2021-10-05 13:39:28.486277",
    "2021-10-05 13:39:28.486281",
    "2021-10-05 13:39:28.486284",
    "2021-10-05 13:39:28.486287",
    "2021-10-05 13:39:28.486290",
    "2021-10-05 13:39:28.486293",
    "2021-10-05 13:39:28.486296",
    "2021-10-05 13:39:28.486299",
    "2021-10-05 13:39:28.486302",
    "2021-10-05 13:39:28.486305",
    "2021-10-05 13:39:28.486308",
    "2021-10-05 13:39:28.486311",
    "2021-10-05 13:39:28.486315",
    "2021-10-05 13:39:28.486318",
    "2021-10-05 13:39:28.486321",
    "2021-10-05 13:39:28.486324",
    "2021-10-05 13:39:28.486327",
    "2021-10-05 13:39:28.486330",
    "2021-10-05 13:39:28.486334",
    "2021-10-05 13:39:28.486337",
    "2021-10-05 13:39:28.486340",
    "2021-10-05 13:39:28.486343",
    "2021-10-05 13:39:28.486347",
    "2021-10-05 13:39:28.486350",
    "2021-10-05 13:39:28.486353",
    "2021-10-05 13:39:28.486356",
    "2021-10-05 13:39:28.486359",
    "2021-10-05 13:39:28.486362",
    "2021-10-05 13:39:28.486365",
    "2021-10-05 13:39:28.486369",
    "2021-10-05 13:39:28.486372",
    "2021-10-05 13:39:28.486375",
    "2021-10-05 13:39:28.486378",
    "2021-10-05 13:39:28.486381",
    "2021-10-05 13:39:28.486384"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
